# Generate Report for Handoff
# Updates the localization-status workbook:
#  - drops the "ffffa99612f2-941e-44be-9143-26876ac55760.md" row
#  - renames 0db4b02f-32a6-4080-8f9f-657f4213dcc2 -> 1fd992ef-0d78-444e-91d1-11892928f7ca
#  - renames c7fea51f-88cc-4160-95a1-e3ebb256e1f3 -> e8ad3b21-98a2-4ea9-b39f-d2e165ae807c
#  - updates handoff file hashes / timestamps

$wb = $excel.ActiveWorkbook

# old 0db4b02f-32a6-4080-8f9f-657f4213dcc2 -> new 1fd992ef-0d78-444e-91d1-11892928f7ca
$newUuid1 = "1fd992ef-0d78-444e-91d1-11892928f7ca"
# old c7fea51f-88cc-4160-95a1-e3ebb256e1f3 -> new e8ad3b21-98a2-4ea9-b39f-d2e165ae807c
$newUuid2 = "e8ad3b21-98a2-4ea9-b39f-d2e165ae807c"

$newHash1 = "f41338ed9099d1b7edce9932d0ca9883d3edc2ac"
$newHash2 = "1cd83df85a2b604f745679ab106928d8ca79f1b0"

$newDateZhCn = "2016-02-29 13:44:23"
$newDateDeDe = "2016-02-29 13:44:35"

$mdName1 = "$newUuid1.md"
$mdName2 = "$newUuid2.md"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/3568cf74f1d5dd8ca3258010779d7773a09e004c/e2e/$mdName1"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/7449d374c63bea7659857d0ae5af8e022ddef9b5/e2e/$mdName2"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7449d374c63bea7659857d0ae5af8e022ddef9b5/.localization-config"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the ffffa99612f2 row (old row 4); row 5 (.localization-config) shifts up to row 4
$ws1.Rows.Item(4).Delete()

# Update file name cells for the two remaining tracked files
$ws1.Range("A2").Value = $mdName1
$ws1.Range("A3").Value = $mdName2

# Rebuild hyperlinks to match the new file names / targets
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrl1, "", "", $mdName1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrl2, "", "", $mdName2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(4).Delete()

$ws2.Range("A2").Value = $mdName1
$ws2.Range("A3").Value = $mdName2

$xlf1ZhCn = "$newUuid1.$newHash1.zh-cn.xlf"
$xlf2ZhCn = "$newUuid2.$newHash2.zh-cn.xlf"

$ws2.Range("C2").Value = $xlf1ZhCn
$ws2.Range("D2").Value = $newDateZhCn
$ws2.Range("C3").Value = $xlf2ZhCn
$ws2.Range("D3").Value = $newDateZhCn

$xlf1ZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ea103d8d354fbc3a87974e227645a3997df911a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlf1ZhCn"
$xlf2ZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8ab69de14e757227a5043455bf570bb6870ec403/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlf2ZhCn"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl1, "", "", $mdName1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $xlf1ZhCnUrl, "", "", $xlf1ZhCn)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrl2, "", "", $mdName2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $xlf2ZhCnUrl, "", "", $xlf2ZhCn)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(4).Delete()

$ws3.Range("A2").Value = $mdName1
$ws3.Range("A3").Value = $mdName2

$xlf1DeDe = "$newUuid1.$newHash1.de-de.xlf"
$xlf2DeDe = "$newUuid2.$newHash2.de-de.xlf"

$ws3.Range("C2").Value = $xlf1DeDe
$ws3.Range("D2").Value = $newDateDeDe
$ws3.Range("C3").Value = $xlf2DeDe
$ws3.Range("D3").Value = $newDateDeDe

$xlf1DeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b5d75a1bb342bf58511e565e175b77d9c41040b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlf1DeDe"
$xlf2DeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95e34de48db4173efcdf38bc1077808aca6ce82c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlf2DeDe"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl1, "", "", $mdName1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $xlf1DeDeUrl, "", "", $xlf1DeDe)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrl2, "", "", $mdName2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $xlf2DeDeUrl, "", "", $xlf2DeDe)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", ".localization-config")
